# Merge the three separate runs describing the "vetor de caracteres" sentence
# into a single run by replacing the whole sentence text via Find/Replace.
$d = $word.ActiveDocument
$d.Content.Find.Execute(
    "Para representar uma string em C, devemos criar um vetor de caracteres, ou seja, um vetor do tipo char.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Para representar uma string em C, devemos criar um vetor de caracteres, ou seja, um vetor do tipo char.",
    2) | Out-Null
